# Update the "Montant adhésion année n" (M), "Montant adhésion année n+1" (N),
# "Montant dons année n" (O) and "Total année n" (P) columns for rows 2-25
# on the active sheet, according to the computed membership/donation amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2;  M=18; N=18; O=124; P=160 },
    @{ Row=3;  M=18; N=18; O=14;  P=50  },
    @{ Row=4;  M=18; N=18; O=164; P=200 },
    @{ Row=5;  M=18; N=18; O=264; P=300 },
    @{ Row=6;  M=18; N=18; O=114; P=150 },
    @{ Row=7;  M=18; N=18; O=14;  P=50  },
    @{ Row=8;  M=18; N=0;  O=12;  P=30  },
    @{ Row=9;  M=18; N=0;  O=12;  P=30  },
    @{ Row=10; M=18; N=18; O=19;  P=55  },
    @{ Row=11; M=18; N=18; O=36;  P=72  },
    @{ Row=12; M=18; N=0;  O=12;  P=30  },
    @{ Row=13; M=18; N=18; O=114; P=150 },
    @{ Row=14; M=18; N=18; O=44;  P=80  },
    @{ Row=15; M=18; N=0;  O=2;   P=20  },
    @{ Row=16; M=18; N=18; O=64;  P=100 },
    @{ Row=17; M=18; N=18; O=64;  P=100 },
    @{ Row=18; M=18; N=18; O=114; P=150 },
    @{ Row=19; M=18; N=18; O=114; P=150 },
    @{ Row=20; M=18; N=18; O=264; P=300 },
    @{ Row=21; M=18; N=18; O=264; P=300 },
    @{ Row=22; M=18; N=18; O=14;  P=50  },
    @{ Row=23; M=18; N=18; O=64;  P=100 },
    @{ Row=24; M=18; N=18; O=64;  P=100 },
    @{ Row=25; M=18; N=18; O=114; P=150 }
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Range("M$row").Value = $u.M
    $ws.Range("N$row").Value = $u.N
    $ws.Range("O$row").Value = $u.O
    $ws.Range("P$row").Value = $u.P
}
